$d = $word.ActiveDocument

# Merge "liderança, " + "coaching" runs (drop the spell-check proofErr split)
$d.Content.Find.Execute("liderança, coaching", $true, $false, $false, $false, $false,
                         $true, 1, $false, "liderança, coaching", 2)

# Merge "vendas centrada" + " no cliente" runs into the preceding run (drop the
# grammar-check proofErr split)
$d.Content.Find.Execute("Total conhecimento do processo de vendas centrada no cliente", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Total conhecimento do processo de vendas centrada no cliente", 2)

# Merge "Defesa dos produtos de " + "georreferenciamento" + " para área agrícola e
# logística" runs (drop the spell-check proofErr split)
$d.Content.Find.Execute("Defesa dos produtos de georreferenciamento para área agrícola e logística", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Defesa dos produtos de georreferenciamento para área agrícola e logística", 2)
